# Insert a new daily price record as row 84 ("Vega Modelo de Temuco" / Puerro),
# which pushes the existing rows 84-216 down to 85-217.
#
# This mirrors the diff: the sheet's used range grows from A1:R216 to A1:R217,
# and every row from the old row 84 onward is shifted down by one position,
# with a brand-new record occupying the vacated row 84.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 84..216 down to 85..217 by inserting a blank row at 84.
$ws.Rows("84").Insert()

# Populate the newly inserted row 84 with the new record's data.
$ws.Range("A84").Value2 = 10
$ws.Range("B84").Value2 = "Vega Modelo de Temuco"
$ws.Range("C84").Value2 = "La Araucanía"
$ws.Range("D84").Value2 = 44757
$ws.Range("E84").Value2 = 9
$ws.Range("F84").Value2 = 100112005
$ws.Range("G84").Value2 = "Puerro"
$ws.Range("H84").Value2 = "Azul de Maquehue"
$ws.Range("I84").Value2 = "Primera"
$ws.Range("J84").Value2 = 20
$ws.Range("K84").Value2 = 17000
$ws.Range("L84").Value2 = 17000
$ws.Range("M84").Value2 = 17000
$ws.Range("N84").Value2 = "$/docena de paquetes"
$ws.Range("O84").Value2 = "Provincia de Cautín"
$ws.Range("P84").Value2 = 1417
$ws.Range("Q84").Value2 = 12
$ws.Range("R84").Value2 = "Hortaliza"

# Ensure the date cell keeps the same date/time number format used by the
# rest of column D (it is copied automatically by Insert, but set it
# explicitly to be safe).
$ws.Range("D84").NumberFormat = $ws.Range("D85").NumberFormat
